{"js": "// Replace the date line and every \"a\u00d7b=c\" multiplication answer in the\n// table with the new values from the target revision. Every old value is\n// unique within the document, so an exact (case-sensitive, non-wildcard)\n// search-and-replace is unambiguous and safe.\nconst replacements = [\n  [\"2024-12-27 Friday\", \"2024-12-28 Saturday\"],\n  [\"46\u00d760=2760\", \"82\u00d797=7954\"],\n  [\"70\u00d793=6510\", \"44\u00d777=3388\"],\n  [\"43\u00d780=3440\", \"52\u00d758=3016\"],\n  [\"88\u00d790=7920\", \"47\u00d750=2350\"],\n  [\"89\u00d766=5874\", \"54\u00d797=5238\"],\n  [\"94\u00d740=3760\", \"84\u00d785=7140\"],\n  [\"67\u00d799=6633\", \"88\u00d741=3608\"],\n  [\"32\u00d721=672\", \"42\u00d771=2982\"],\n  [\"13\u00d714=182\", \"19\u00d781=1539\"],\n  [\"99\u00d747=4653\", \"86\u00d778=6708\"],\n  [\"81\u00d799=8019\", \"48\u00d731=1488\"],\n  [\"83\u00d711=913\", \"30\u00d795=2850\"],\n  [\"33\u00d749=1617\", \"24\u00d744=1056\"],\n  [\"79\u00d719=1501\", \"79\u00d777=6083\"],\n  [\"57\u00d725=1425\", \"33\u00d776=2508\"],\n  [\"80\u00d792=7360\", \"90\u00d777=6930\"],\n  [\"94\u00d767=6298\", \"74\u00d729=2146\"],\n  [\"72\u00d733=2376\", \"66\u00d725=1650\"],\n  [\"91\u00d798=8918\", \"59\u00d758=3422\"],\n  [\"69\u00d769=4761\", \"25\u00d749=1225\"],\n  [\"20\u00d761=1220\", \"20\u00d755=1100\"],\n  [\"26\u00d761=1586\", \"76\u00d776=5776\"],\n  [\"51\u00d776=3876\", \"38\u00d782=3116\"],\n  [\"73\u00d726=1898\", \"46\u00d732=1472\"],\n  [\"42\u00d751=2142\", \"84\u00d759=4956\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and every \"a\u00d7b=c\" multiplication answer in the\n# table with the new values from the target revision. Every old value is\n# unique within the document, so an exact (case-sensitive, non-wildcard)\n# Find/Replace is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-12-27 Friday\", \"2024-12-28 Saturday\"),\n    @(\"46\u00d760=2760\", \"82\u00d797=7954\"),\n    @(\"70\u00d793=6510\", \"44\u00d777=3388\"),\n    @(\"43\u00d780=3440\", \"52\u00d758=3016\"),\n    @(\"88\u00d790=7920\", \"47\u00d750=2350\"),\n    @(\"89\u00d766=5874\", \"54\u00d797=5238\"),\n    @(\"94\u00d740=3760\", \"84\u00d785=7140\"),\n    @(\"67\u00d799=6633\", \"88\u00d741=3608\"),\n    @(\"32\u00d721=672\", \"42\u00d771=2982\"),\n    @(\"13\u00d714=182\", \"19\u00d781=1539\"),\n    @(\"99\u00d747=4653\", \"86\u00d778=6708\"),\n    @(\"81\u00d799=8019\", \"48\u00d731=1488\"),\n    @(\"83\u00d711=913\", \"30\u00d795=2850\"),\n    @(\"33\u00d749=1617\", \"24\u00d744=1056\"),\n    @(\"79\u00d719=1501\", \"79\u00d777=6083\"),\n    @(\"57\u00d725=1425\", \"33\u00d776=2508\"),\n    @(\"80\u00d792=7360\", \"90\u00d777=6930\"),\n    @(\"94\u00d767=6298\", \"74\u00d729=2146\"),\n    @(\"72\u00d733=2376\", \"66\u00d725=1650\"),\n    @(\"91\u00d798=8918\", \"59\u00d758=3422\"),\n    @(\"69\u00d769=4761\", \"25\u00d749=1225\"),\n    @(\"20\u00d761=1220\", \"20\u00d755=1100\"),\n    @(\"26\u00d761=1586\", \"76\u00d776=5776\"),\n    @(\"51\u00d776=3876\", \"38\u00d782=3116\"),\n    @(\"73\u00d726=1898\", \"46\u00d732=1472\"),\n    @(\"42\u00d751=2142\", \"84\u00d759=4956\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
